# Insert a new data row at row 238 (pushing existing rows 238:265 down to 239:266)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(238).Insert()

$ws.Cells.Item(238, 1).Value = 5
$ws.Cells.Item(238, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(238, 3).Value = "Maule"
$ws.Cells.Item(238, 4).Value = 44748
$ws.Cells.Item(238, 5).Value = 7
$ws.Cells.Item(238, 6).Value = "Fruta"
$ws.Cells.Item(238, 7).Value = 100108
$ws.Cells.Item(238, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(238, 9).Value = 100108005
$ws.Cells.Item(238, 10).Value = "Piña"
$ws.Cells.Item(238, 11).Value = "Caramelo"
$ws.Cells.Item(238, 12).Value = "Segunda"
$ws.Cells.Item(238, 13).Value = 200
$ws.Cells.Item(238, 14).Value = 19000
$ws.Cells.Item(238, 15).Value = 19000
$ws.Cells.Item(238, 16).Value = 19000
$ws.Cells.Item(238, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(238, 18).Value = "Ecuador"
$ws.Cells.Item(238, 19).Value = 1357
$ws.Cells.Item(238, 20).Value = 14
